$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.953.84'
$ws.Range('E2').Value = '  +1.37%  '
$ws.Range('D3').Value = '1.767.89'
$ws.Range('E3').Value = '  +0.86%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '329.10'
$ws.Range('E6').Value = '  -0.07%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4560'
$ws.Range('E7').Value = '  +1.44%  '
$ws.Range('E8').Value = '  -1.13%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '42.10'
$ws.Range('E9').Value = '  +1.62%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07410'
$ws.Range('E10').Value = '  -0.81%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.098'
$ws.Range('E11').Value = '  +1.42%  '
$ws.Range('E12').Value = '  -0.12%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.75'
$ws.Range('E13').Value = '  +0.06%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.004'
$ws.Range('E14').Value = '  +0.28%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.194'
$ws.Range('E15').Value = '  +0.60%  '
$ws.Range('D16').Value = '1.777.63'
$ws.Range('E16').Value = '  +1.27%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '92.62'
$ws.Range('E17').Value = '  -1.13%  '
$ws.Range('E18').Value = '  +0.30%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06452'
$ws.Range('E19').Value = '  +1.14%  '
$ws.Range('E20').Value = '  -0.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.94'
$ws.Range('E21').Value = '  -1.22%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.777'
$ws.Range('D23').Value = '27.977.38'
$ws.Range('E23').Value = '  +1.28%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.23'
$ws.Range('E24').Value = '  +0.27%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.135'
$ws.Range('E25').Value = '  +2.34%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '161.87'
$ws.Range('E26').Value = '  -2.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.20'
$ws.Range('E27').Value = '  +0.12%  '
$ws.Range('D28').Value = '1.973.91'
$ws.Range('E28').Value = '  +0.91%  '
$ws.Range('E29').Value = '  +2.97%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '123.97'
$ws.Range('E30').Value = '  -1.19%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.080'
$ws.Range('E31').Value = '  -1.06%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09305'
$ws.Range('E32').Value = '  +1.26%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.605'
$ws.Range('E33').Value = '  +1.79%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.670'
$ws.Range('E34').Value = '  +0.41%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '11.87'
$ws.Range('E35').Value = '  +0.70%  '
$ws.Range('E36').Value = '  -0.42%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06138'
$ws.Range('E37').Value = '  +2.01%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2086'
$ws.Range('E38').Value = '  -0.14%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.952'
$ws.Range('E39').Value = '  +0.37%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6274'
$ws.Range('E40').Value = '  -0.25%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.185'
$ws.Range('E41').Value = '  +0.18%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.381'
$ws.Range('E42').Value = '  -1.37%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '7.822'
$ws.Range('E43').Value = '  +0.33%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.22'
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('E45').Value = '  +0.51%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5875'
$ws.Range('E46').Value = '  +0.11%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '122.59'
$ws.Range('E47').Value = '  +0.45%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.941'
$ws.Range('E48').Value = '  +0.33%  '
$ws.Range('E49').Value = '  -0.11%  '
$ws.Range('E50').Value = '  -0.90%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '73.32'
$ws.Range('E51').Value = '  +2.37%  '
